$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" '34.546.32'
$ws.Range("E2").Value = '  -0.01%  '
Set-TextCell "D3" '1.809.81'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.11%  '
Set-TextCell "D5" '225.91'
$ws.Range("E5").Value = '  -1.43%  '
Set-TextCell "D6" '0.598'
$ws.Range("E6").Value = '  +3.11%  '
$ws.Range("E7").Value = '  -0.10%  '
Set-TextCell "D8" '37.49'
$ws.Range("E8").Value = '  +7.37%  '
$ws.Range("E9").Value = '  -3.46%  '
Set-TextCell "D10" '0.0680'
$ws.Range("E10").Value = '  -2.34%  '
Set-TextCell "D11" '0.0968'
$ws.Range("E11").Value = '  +1.23%  '
Set-TextCell "D12" '2.070.12'
$ws.Range("E12").Value = '  -0.21%  '
Set-TextCell "D13" '11.36'
$ws.Range("E13").Value = '  +1.29%  '
Set-TextCell "D14" '1.824.43'
$ws.Range("E14").Value = '  +0.69%  '
Set-TextCell "D15" '0.634'
$ws.Range("E15").Value = '  -2.46%  '
Set-TextCell "D16" '34.529.70'
$ws.Range("E16").Value = '  -0.06%  '
Set-TextCell "D17" '4.44'
$ws.Range("E17").Value = '  -0.83%  '
Set-TextCell "D18" '68.67'
$ws.Range("E18").Value = '  -0.87%  '
Set-TextCell "D19" '243.77'
$ws.Range("E19").Value = '  -0.86%  '
Set-TextCell "D20" '0.0₃0776'
$ws.Range("E20").Value = '  -2.98%  '
Set-TextCell "D21" '11.24'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("E22").Value = '  -0.09%  '
Set-TextCell "D23" '4.14'
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("E24").Value = '  +3.85%  '
Set-TextCell "D25" '171.93'
$ws.Range("E25").Value = '  -0.55%  '
Set-TextCell "D26" '7.84'
$ws.Range("E26").Value = '  -1.56%  '
$ws.Range("E27").Value = '  +2.69%  '
Set-TextCell "D28" '0.121'
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell "D30" '3.94'
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell "D31" '3.82'
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("E32").Value = '  -1.41%  '
Set-TextCell "D33" '0.0520'
$ws.Range("E33").Value = '  -2.74%  '
Set-TextCell "D34" '1.83'
$ws.Range("E34").Value = '  -1.34%  '
Set-TextCell "D35" '1.366.02'
$ws.Range("E35").Value = '  -2.32%  '
Set-TextCell "D36" '0.653'
$ws.Range("E36").Value = '  -4.59%  '
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("E38").Value = '  -4.65%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell "D40" '1.21'
$ws.Range("E40").Value = '  +8.09%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell "D41" '2.43'
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D42" '80.93'
$ws.Range("E42").Value = '  -3.46%  '
Set-TextCell "D43" '0.941'
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell "D44" '2.77'
$ws.Range("E44").Value = '  -2.10%  '
Set-TextCell "D45" '13.87'
$ws.Range("E45").Value = '  +2.98%  '
$ws.Range("E46").Value = '  -2.89%  '
Set-TextCell "D47" '1.971.00'
$ws.Range("E47").Value = '  -0.17%  '
Set-TextCell "D48" '5.83'
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("E49").Value = '  -0.11%  '
Set-TextCell "D50" '103.01'
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell "D51" '0.0₆0122'
$ws.Range("E51").Value = '  -6.93%  '
